# Apply weekly fruta/hortaliza price update: swap data between row pairs
# (2 <-> 7), (3 <-> 6), (5 <-> 9) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($Worksheet, $Col, $RowA, $RowB) {
    $cellA = $Worksheet.Range("$Col$RowA")
    $cellB = $Worksheet.Range("$Col$RowB")

    $valueA = $cellA.Value2
    $valueB = $cellB.Value2

    $cellA.Value2 = $valueB
    $cellB.Value2 = $valueA
}

# Row 2 <-> Row 7: Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg
foreach ($col in @("D", "M", "N", "O", "P", "S")) {
    Swap-Cell $ws $col 2 7
}

# Row 3 <-> Row 6: Fecha, Volumen
foreach ($col in @("D", "M")) {
    Swap-Cell $ws $col 3 6
}

# Row 5 <-> Row 9: Fecha, Volumen
foreach ($col in @("D", "M")) {
    Swap-Cell $ws $col 5 9
}
